$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats
$xlPasteFormats = -4122

# ---------------------------------------------------------------------------
# Row 20 used to be the last data row (a "continuation" line for filename
# 307, reusing the bottom-border style of the block). Now that a new row is
# appended below it, row 20 must be restyled from a "last row of block"
# style (s=4/5) into a "continuation row with bottom border" style (s=6/7),
# matching the look of rows 3/5 (the other continuation rows in the sheet).
# Copy formats only, so the existing values/shared-string refs are kept.
# ---------------------------------------------------------------------------
$ws.Range("A3").Copy()
$ws.Range("A20").PasteSpecial($xlPasteFormats)
$ws.Range("B3").Copy()
$ws.Range("B20").PasteSpecial($xlPasteFormats)
$ws.Range("C3:E3").Copy()
$ws.Range("C20:E20").PasteSpecial($xlPasteFormats)

# ---------------------------------------------------------------------------
# Row 21: brand-new table entry.
#   A = filename, B = line number, C = English, D = translated (RU),
#   E = "converted" (cipher) string.
# ---------------------------------------------------------------------------
$ws.Range("C21").Value = " Yes, [CS:N]Loudred[CR]!"
$ws.Range("A21").Value = "SCRIPT/G01P04A/um1604.ssb"
$ws.Range("D21").Value = " Да, [CS:N]Лаудред[CR]!"
$ws.Range("E21").Value = " Äà, [CS:N]Ìàôäñåä[CR]!"
$ws.Range("B21").Value = 252

# Give row 21 the normal "first/only row of a block" style (s=4/5), matching
# row 2 (the first such block in the sheet).
$ws.Range("A2").Copy()
$ws.Range("A21").PasteSpecial($xlPasteFormats)
$ws.Range("B2").Copy()
$ws.Range("B21").PasteSpecial($xlPasteFormats)
$ws.Range("C2:E2").Copy()
$ws.Range("C21:E21").PasteSpecial($xlPasteFormats)

# Row 21 wraps onto three lines at this column width -> taller row.
$ws.Rows.Item(21).RowHeight = 43.2

# ---------------------------------------------------------------------------
# View state: scroll down so the new rows are visible, and move the
# selection onto the newly-added block, like the original author did.
# ---------------------------------------------------------------------------
$win = $ws.Application.ActiveWindow
$win.ScrollRow = 16
$win.ScrollColumn = 1
$ws.Range("C19").Select() | Out-Null
